# Auto-generated: bulk numeric value refresh across the Ragnarok Profits
# workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), mirroring the
# "update Sheets via scheduled runner" data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 9
$ws.Range("H9").Value = 271.35715
$ws.Range("I9").Value = 267.6
$ws.Range("J9").Value = 280.75
$ws.Range("K9").Value = 267.6
$ws.Range("L9").Value = 280.75
$ws.Range("M9").Value = -98.60000000000002
$ws.Range("N9").Value = -618.75

# ALC row 19
$ws.Range("H19").Value = 22808410
$ws.Range("J19").Value = 8334995.5
$ws.Range("L19").Value = 8334995.5
$ws.Range("N19").Value = -8335345.5

# ALC row 21
$ws.Range("H21").Value = 40999.5
$ws.Range("I21").Value = 40999.5
$ws.Range("K21").Value = 40999.5
$ws.Range("M21").Value = -40531.5

# ALC row 23
$ws.Range("H23").Value = 40999.5
$ws.Range("I23").Value = 40999.5
$ws.Range("K23").Value = 40999.5
$ws.Range("M23").Value = -40765.5

# ALC row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# ALC row 40
$ws.Range("H40").Value = 50002390
$ws.Range("J40").Value = 125002990
$ws.Range("L40").Value = 125002990
$ws.Range("N40").Value = -125003340

# ALC row 94
$ws.Range("H94").Value = 7265.8
$ws.Range("I94").Value = 7265.8
$ws.Range("K94").Value = 7265.8
$ws.Range("M94").Value = -6814.8

# ALC row 106
$ws.Range("H106").Value = 11353.182
$ws.Range("I106").Value = 11235.625
$ws.Range("K106").Value = 11235.625
$ws.Range("M106").Value = -10604.625

# ALC row 132
$ws.Range("H132").Value = 4272.2617
$ws.Range("I132").Value = 1935.931
$ws.Range("J132").Value = 9484.076999999999
$ws.Range("K132").Value = 5807.793
$ws.Range("L132").Value = 28452.231
$ws.Range("M132").Value = -3277.793
$ws.Range("N132").Value = -33512.231

# ALC row 137
$ws.Range("H137").Value = 1182630.8
$ws.Range("I137").Value = 2188
$ws.Range("K137").Value = 6564
$ws.Range("M137").Value = -4014

$ws = $wb.Worksheets.Item("ARM")
# ARM row 16
$ws.Range("H16").Value = 1633.3334
$ws.Range("I16").Value = 1550
$ws.Range("J16").Value = 1675
$ws.Range("K16").Value = 1550
$ws.Range("L16").Value = 1675
$ws.Range("M16").Value = -1263
$ws.Range("N16").Value = -2249

# ARM row 45
$ws.Range("H45").Value = 1526.1818
$ws.Range("I45").Value = 1478.9
$ws.Range("K45").Value = 1478.9
$ws.Range("M45").Value = -1101.9

# ARM row 122
$ws.Range("H122").Value = 3639
$ws.Range("J122").Value = 3854
$ws.Range("L122").Value = 11562
$ws.Range("N122").Value = -16462

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 7334.067
$ws.Range("I20").Value = 10240.777
$ws.Range("K20").Value = 10240.777
$ws.Range("M20").Value = -9993.777

# BSM row 80
$ws.Range("H80").Value = 860.8570999999999
$ws.Range("I80").Value = 540
$ws.Range("K80").Value = 540
$ws.Range("M80").Value = 458

# BSM row 83
$ws.Range("H83").Value = 860.8570999999999
$ws.Range("I83").Value = 540
$ws.Range("K83").Value = 2700
$ws.Range("M83").Value = 2292

# BSM row 94
$ws.Range("H94").Value = 3554.05
$ws.Range("I94").Value = 3438.7334
$ws.Range("J94").Value = 3900
$ws.Range("K94").Value = 3438.7334
$ws.Range("L94").Value = 3900
$ws.Range("M94").Value = -2987.7334
$ws.Range("N94").Value = -4802

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 168354770
$ws.Range("I31").Value = 200005520
$ws.Range("J31").Value = 10101010
$ws.Range("K31").Value = 200005520
$ws.Range("L31").Value = 10101010
$ws.Range("M31").Value = -200005225
$ws.Range("N31").Value = -10101600

# CRP row 34
$ws.Range("H34").Value = 168354770
$ws.Range("I34").Value = 200005520
$ws.Range("J34").Value = 10101010
$ws.Range("K34").Value = 200005520
$ws.Range("L34").Value = 10101010
$ws.Range("M34").Value = -200005318
$ws.Range("N34").Value = -10101414

# CRP row 99
$ws.Range("H99").Value = 21617
$ws.Range("I99").Value = 11927.286
$ws.Range("K99").Value = 11927.286
$ws.Range("M99").Value = -10429.286

# CRP row 107
$ws.Range("H107").Value = 4000.577
$ws.Range("J107").Value = 4506.222
$ws.Range("L107").Value = 4506.222
$ws.Range("N107").Value = -8346.222

# CRP row 126
$ws.Range("H126").Value = 21617
$ws.Range("I126").Value = 11927.286
$ws.Range("K126").Value = 35781.858
$ws.Range("M126").Value = -33311.858

$ws = $wb.Worksheets.Item("CUL")
# CUL row 2
$ws.Range("H2").Value = 136.09091

# CUL row 38
$ws.Range("H38").Value = 165.66667
$ws.Range("I38").Value = 8
$ws.Range("K38").Value = 24
$ws.Range("M38").Value = 323

$ws = $wb.Worksheets.Item("GSM")
# GSM row 123
$ws.Range("H123").Value = 99999.336
$ws.Range("J123").Value = 99999.336
$ws.Range("L123").Value = 99999.336
$ws.Range("N123").Value = -104899.336

# GSM row 126
$ws.Range("H126").Value = 2038.9412
$ws.Range("I126").Value = 1697.625
$ws.Range("K126").Value = 5092.875
$ws.Range("M126").Value = -2622.875

# GSM row 132
$ws.Range("H132").Value = 4548697
$ws.Range("I132").Value = 3034.4285
$ws.Range("K132").Value = 9103.2855
$ws.Range("M132").Value = -6573.2855

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 3993.95
$ws.Range("I16").Value = 2406.75
$ws.Range("J16").Value = 6374.75
$ws.Range("K16").Value = 2406.75
$ws.Range("L16").Value = 6374.75
$ws.Range("M16").Value = -2236.75
$ws.Range("N16").Value = -6714.75

# LTW row 42
$ws.Range("H42").Value = 25555.445
$ws.Range("J42").Value = 18571.428
$ws.Range("L42").Value = 18571.428
$ws.Range("N42").Value = -19697.428

# LTW row 49
$ws.Range("H49").Value = 25555.445
$ws.Range("J49").Value = 18571.428
$ws.Range("L49").Value = 18571.428
$ws.Range("N49").Value = -18865.428

$ws = $wb.Worksheets.Item("WVR")
# WVR row 38
$ws.Range("H38").Value = 90000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 90000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 90000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -90946

# WVR row 48
$ws.Range("H48").Value = 89999
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 89999
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 89999
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -91137

# WVR row 49
$ws.Range("H49").Value = 125000
$ws.Range("J49").Value = 125000
$ws.Range("L49").Value = 125000
$ws.Range("N49").Value = -125460

# WVR row 126
$ws.Range("H126").Value = 6145.5264
$ws.Range("I126").Value = 6178.1875
$ws.Range("K126").Value = 18534.5625
$ws.Range("M126").Value = -16064.5625

# WVR row 136
$ws.Range("H136").Value = 229320.44
$ws.Range("I136").Value = 1901.25
$ws.Range("K136").Value = 5703.75
$ws.Range("M136").Value = -3153.75

